# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Column D ("Price") holds free-form text (some values use a '.' as a
# thousands separator, e.g. "26.538.92"), so every D-cell write forces the
# cell to Text format first and restores the default "Normal" style
# afterwards -- this keeps the literal digits/zeros intact (Excel would
# otherwise silently reinterpret "1.000" as the number 1, "0.3220" as
# 0.322, etc.) without leaving a lingering custom number format on the
# cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '26.538.92'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  +0.19%  '

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '1.847.22'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  -0.15%  '

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '1.000'
$cell.Style = "Normal"
$ws.Range("E4").Value = '  -0.06%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '264.04'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +1.36%  '

$ws.Range("E6").Value = '  -0.03%  '

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.5228'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  +1.23%  '

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.3220'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  -0.59%  '

$ws.Range("E9").Value = '  +0.47%  '

$ws.Range("E10").Value = '  -1.02%  '

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.7801'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  +1.05%  '

$ws.Range("E12").Value = '  +0.53%  '

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '1.834.87'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  -1.16%  '

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '88.51'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  -0.18%  '

$ws.Range("E15").Value = '  -0.30%  '

$ws.Range("E16").Value = '  -0.05%  '

$ws.Range("E17").Value = '  -0.85%  '

$ws.Range("B18").Value = 'Dai'
$ws.Range("C18").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '1.000'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  +0.00%  '

$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '0.000007963'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +0.61%  '

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '26.554.32'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  +0.04%  '

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '4.631'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  +2.13%  '

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '9.469'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  -0.66%  '

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '6.005'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  +1.29%  '

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '142.93'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  -1.08%  '

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '2.176'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  -7.36%  '

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '1.689'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  +2.21%  '

$ws.Range("E27").Value = '  +0.11%  '

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '111.69'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  +0.36%  '

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '4.182'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  -0.67%  '

$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '4.118'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  -1.34%  '

$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '0.08737'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  -0.10%  '

$ws.Range("E32").Value = '  +0.46%  '

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '0.7215'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +4.90%  '

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '1.131'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  -0.30%  '

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '2.867'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  +0.68%  '

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '3.108'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  -0.77%  '

$ws.Range("E37").Value = '  -0.15%  '

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '2.220'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  +0.37%  '

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '0.4865'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  -0.82%  '

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '0.8970'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  -0.47%  '

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '110.97'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  -1.72%  '

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '6.028'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  -2.05%  '

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '1.000'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  +0.01%  '

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '7.640'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  -1.87%  '

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '0.4210'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  +0.03%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '9.074'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  -0.19%  '

$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '0.05889'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  -0.02%  '

$ws.Range("E48").Value = '  -1.58%  '

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '35.02'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  -0.75%  '

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '0.8890'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  +3.51%  '

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '59.91'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  +0.96%  '

Write-Output "cryptos list updated"
